$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "opis"
$ws.Range("E2").Value = "opis opis"
$ws.Range("E3").Value = "ccoewniorvewmcsof"

$ws.Range("E3").Select()
